$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 11:22"

# Row 13
$ws.Range("B13").Value = 28018
$ws.Range("C13").Value = 1351
$ws.Range("D13").Value = 5986
$ws.Range("E13").Value = 18686
$ws.Range("F13").Value = 1262
$ws.Range("G13").Value = 327
$ws.Range("H13").Value = 3346

# Row 19
$ws.Range("B19").Value = 13713
$ws.Range("C19").Value = 153
$ws.Range("D19").Value = 6604
$ws.Range("E19").Value = 6772
$ws.Range("F19").Value = 246
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 337

# Row 37
$ws.Range("A37").Value = "Malasia"
$ws.Range("B37").Value = 4530
$ws.Range("C37").Value = 184
$ws.Range("D37").Value = 1995
$ws.Range("E37").Value = 2462
$ws.Range("F37").Value = 72
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 73

# Row 38
$ws.Range("A38").Value = "Filipinas"
$ws.Range("B38").Value = 4428
$ws.Range("C38").Value = 233
$ws.Range("D38").Value = 157
$ws.Range("E38").Value = 4024
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 26
$ws.Range("H38").Value = 247

# Row 46
$ws.Range("B46").Value = 2905
$ws.Range("C46").Value = 136
$ws.Range("D46").Value = 300
$ws.Range("E46").Value = 2557
$ws.Range("F46").Value = 82
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 48

# Row 52
$ws.Range("A52").Value = "Bielorrusia"
$ws.Range("B52").Value = 2226
$ws.Range("C52").Value = 245
$ws.Range("D52").Value = 172
$ws.Range("E52").Value = 2031
$ws.Range("F52").Value = 72
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 23

# Row 53
$ws.Range("A53").Value = "Singapur"
$ws.Range("B53").Value = 2108
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 492
$ws.Range("E53").Value = 1609
$ws.Range("F53").Value = 29
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 7

# Row 54
$ws.Range("A54").Value = "Grecia"
$ws.Range("B54").Value = 2011
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 269
$ws.Range("E54").Value = 1650
$ws.Range("F54").Value = 77
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 92

# Row 55
$ws.Range("A55").Value = "Sudafrica"
$ws.Range("B55").Value = 2003
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 410
$ws.Range("E55").Value = 1569
$ws.Range("F55").Value = 7
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 24

# Row 67
$ws.Range("B67").Value = 1188
$ws.Range("C67").Value = 28
$ws.Range("D67").Value = 148
$ws.Range("E67").Value = 990
$ws.Range("F67").Value = 37
$ws.Range("G67").Value = 5
$ws.Range("H67").Value = 50

# Row 68
$ws.Range("A68").Value = "Kuwait"
$ws.Range("B68").Value = 1154
$ws.Range("C68").Value = 161
$ws.Range("D68").Value = 133
$ws.Range("E68").Value = 1020
$ws.Range("F68").Value = 27
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 1

# Row 69
$ws.Range("A69").Value = "Lituania"
$ws.Range("B69").Value = 1026
$ws.Range("C69").Value = 27
$ws.Range("D69").Value = 54
$ws.Range("E69").Value = 949
$ws.Range("F69").Value = 14
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 23

# Row 70
$ws.Range("A70").Value = "Hong Kong"
$ws.Range("B70").Value = 1001
$ws.Range("C70").Value = 11
$ws.Range("D70").Value = 336
$ws.Range("E70").Value = 661
$ws.Range("F70").Value = 14
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 4

# Row 71
$ws.Range("A71").Value = "Barein"
$ws.Range("B71").Value = 998
$ws.Range("C71").Value = 73
$ws.Range("D71").Value = 551
$ws.Range("E71").Value = 441
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 6

# Row 84
$ws.Range("B84").Value = 619
$ws.Range("C84").Value = 10
$ws.Range("D84").Value = 76
$ws.Range("E84").Value = 523
$ws.Range("F84").Value = 28
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 20

# Row 89
$ws.Range("A89").Value = "Afganistan"
$ws.Range("B89").Value = 555
$ws.Range("C89").Value = 34
$ws.Range("D89").Value = 32
$ws.Range("E89").Value = 505
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 18

# Row 90
$ws.Range("A90").Value = "Oman"
$ws.Range("B90").Value = 546
$ws.Range("C90").Value = 62
$ws.Range("D90").Value = 109
$ws.Range("E90").Value = 434
$ws.Range("F90").Value = 3
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 3
